{"js": "// Style sidebar nav with club primary color\n// Updates the two \"Sidebar\" to-do bullets to describe using the configured\n// (club) color with a fallback to the default orange, instead of plain\n// \"oranje\" styling.\n\nconst body = context.document.body;\n\n// Locate the two target phrases first (search results become stale once the\n// document is mutated, so resolve both ranges before editing either one).\nconst hoverResults = body.search(\"ook oranje bij hover\", { matchCase: false });\nhoverResults.load(\"items\");\n\nconst selectedPageResults = body.search(\n  \"Geselecteerde pagina verticaal streepje voor de tekst en kleur niet oranje?\",\n  { matchCase: false }\n);\nselectedPageResults.load(\"items\");\n\nawait context.sync();\n\nif (hoverResults.items.length > 0) {\n  hoverResults.items[0].insertText(\n    \"in menu graag in ingestelde kleur (met fallback naar standaard) bij hover\",\n    Word.InsertLocation.replace\n  );\n}\n\nif (selectedPageResults.items.length > 0) {\n  selectedPageResults.items[0].insertText(\n    \"Geselecteerde pagina en verticaal streepje voor de tekst in de ingestelde kleur met fallback.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Style sidebar nav with club primary color\n# Updates the two \"Sidebar\" to-do bullets to describe using the configured\n# (club) color with a fallback to the default orange, instead of plain\n# \"oranje\" styling.\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.Text = \"ook oranje bij hover\"\n$find1.Replacement.Text = \"in menu graag in ingestelde kleur (met fallback naar standaard) bij hover\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.Text = \"Geselecteerde pagina verticaal streepje voor de tekst en kleur niet oranje?\"\n$find2.Replacement.Text = \"Geselecteerde pagina en verticaal streepje voor de tekst in de ingestelde kleur met fallback.\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
